# Add Publons test script row (PUBLONS022 / OPQA-5890) to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("A22").Value = "PUBLONS022"
$ws.Range("B22").Value = "OPQA-5890"
$ws.Range("C22").Value = "Verify Error message when user  exist in platform in suspend state and trying to login."
$ws.Range("D22").Value = "Y"

# Update the view: scroll so row 7 is at the top, and select D22 (matches the
# author's in-progress edit position when the workbook was saved).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D22").Select()
